$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure column D retains text formatting when assigning numeric-looking strings
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '59.096.61'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.310.85'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '541.53'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').Value = '132.51'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('D9').Value = '2.312.61'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('E10').Value = '  -1.20%  '
$ws.Range('D11').Value = '5.43'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '0.332'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = '23.86'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '2.727.91'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '59.012.71'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '2.304.89'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').Value = '10.61'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').Value = '4.17'
$ws.Range('E20').Value = '  -3.30%  '
$ws.Range('D21').Value = '312.81'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').Value = '6.58'
$ws.Range('E22').Value = '  +1.88%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '62.70'
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('D25').Value = '0.172'
$ws.Range('E25').Value = '  +2.39%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '7.94'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').Value = '1.30'
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('D30').Value = '170.04'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('E31').Value = '  +5.92%  '
$ws.Range('D32').Value = '0.0₃0738'
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('D33').Value = '5.90'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').Value = '0.384'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '17.84'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').Value = '1.32'
$ws.Range('E37').Value = '  +4.62%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = '4.06'
$ws.Range('E39').Value = '  +2.19%  '
$ws.Range('D40').Value = '38.42'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').Value = '309.33'
$ws.Range('E41').Value = '  +3.51%  '
$ws.Range('D42').Value = '1.52'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = '140.91'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = '3.45'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = '0.0957'
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('D46').Value = '0.0495'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').Value = '0.557'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').Value = '18.40'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').Value = '0.0211'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '11.00'
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('E51').Value = '  -0.14%  '
